# Update bouquet table: add handling for new bouquets (ורד אדום קטן family),
# delete now-obsolete rows, and refresh the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that are no longer part of the bouquet list (old rows 8,9,10).
$ws.Rows.Item(8).Resize(3).Delete()

# Row 2: בוקט ורדים גדול (based on ורד אדום קטן) / גפסנית / לבן / רגיל / 5
$ws.Range("A2").Value = "בוקט ורדים גדול (based on ורד אדום קטן)"
$ws.Range("B2").Value = "גפסנית"
$ws.Range("C2").Value = "לבן"
$ws.Range("D2").Value = "רגיל"
$ws.Range("E2").Value = 5

# Row 3: בוקט ורדים גדול (based on ורד אדום קטן) / ורד / אדום / גדול / 10
$ws.Range("A3").Value = "בוקט ורדים גדול (based on ורד אדום קטן)"
$ws.Range("B3").Value = "ורד"
$ws.Range("C3").Value = "אדום"
$ws.Range("D3").Value = "גדול"
$ws.Range("E3").Value = 10

# Row 4: ורד אדום קטן / ורד / אדום / קטן / 5
$ws.Range("A4").Value = "ורד אדום קטן"
$ws.Range("B4").Value = "ורד"
$ws.Range("C4").Value = "אדום"
$ws.Range("D4").Value = "קטן"
$ws.Range("E4").Value = 5

# Row 5: ורד אדום קטן / גפסנית / לבן / רגיל / 2
$ws.Range("A5").Value = "ורד אדום קטן"
$ws.Range("B5").Value = "גפסנית"
$ws.Range("C5").Value = "לבן"
$ws.Range("D5").Value = "רגיל"
$ws.Range("E5").Value = 2

# Row 6: ורד ענק / ורד / אדום / קטן / 15
$ws.Range("A6").Value = "ורד ענק"
$ws.Range("B6").Value = "ורד"
$ws.Range("C6").Value = "אדום"
$ws.Range("D6").Value = "קטן"
$ws.Range("E6").Value = 15

# Row 7: ורד ענק / גפסנית / לבן / רגיל / 10
$ws.Range("A7").Value = "ורד ענק"
$ws.Range("B7").Value = "גפסנית"
$ws.Range("C7").Value = "לבן"
$ws.Range("D7").Value = "רגיל"
$ws.Range("E7").Value = 10
